# Lecture partielle de l'EDT M1 MIAGE.
# Shift the sample week from Jan 2023 to the Jan 2026 occurrence that
# preserves the same weekday pattern (Thu, Fri, Mon, Fri, Mon) and update
# the displayed day-name labels in column B to match the new dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dimanche 08/01/2023 -> jeudi 08/01/2026
$ws.Range("A2").Value = 46030
$ws.Range("B2").Value = "jeudi"

# Row 4: lundi 09/01/2023 -> vendredi 09/01/2026
$ws.Range("A4").Value = 46031
$ws.Range("B4").Value = "vendredi"

# Row 7: jeudi 12/01/2023 -> lundi 12/01/2026
$ws.Range("A7").Value = 46034
$ws.Range("B7").Value = "lundi"

# Row 10: lundi 16/01/2023 -> vendredi 16/01/2026
$ws.Range("A10").Value = 46038
$ws.Range("B10").Value = "vendredi"

# Row 13: jeudi 19/01/2023 -> lundi 19/01/2026
$ws.Range("A13").Value = 46041
$ws.Range("B13").Value = "lundi"
